$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("taxonomy_translations")

$rows = @(
    @(14, "en", "name", "Car Brand"),
    @(15, "en", "name", "Car Model"),
    @(16, "en", "name", "Car Generation"),
    @(17, "en", "name", "Car Body Style"),
    @(18, "en", "name", "Car Engine Capacity"),
    @(19, "en", "name", "Powertrain"),
    @(20, "en", "name", "Engine Type"),
    @(21, "en", "name", "Exterior Color"),
    @(22, "en", "name", "Interior Color"),
    @(23, "en", "name", "Horsepower"),
    @(14, "zh_Hant", "name", "廠牌"),
    @(15, "zh_Hant", "name", "車款"),
    @(16, "zh_Hant", "name", "世代"),
    @(17, "zh_Hant", "name", "車身樣式"),
    @(18, "zh_Hant", "name", "排氣量"),
    @(19, "zh_Hant", "name", "動力系統"),
    @(20, "zh_Hant", "name", "引擎類型"),
    @(21, "zh_Hant", "name", "外觀顏色"),
    @(22, "zh_Hant", "name", "內裝顏色"),
    @(23, "zh_Hant", "name", "馬力")
)

$r = 29
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("A33").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.SplitRow = 2
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D55").Select()
